$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear obsolete rows 19-42 (data now fits in rows 1-18)
$ws.Range("A19:F42").Clear()

# Rewrite the A1:F18 block with the recomputed similarity values
$ws.Cells.Item(1, 4).Value = -0.02042428036513438
$ws.Cells.Item(1, 5).Value = -0.2823383113349731
$ws.Cells.Item(1, 6).Value = [double]"-3.338007351689955E-05"

$ws.Cells.Item(2, 1).Value = $null
$ws.Cells.Item(2, 4).Value = -0.04241237699940882
$ws.Cells.Item(2, 5).Value = -0.04748827222003059
$ws.Cells.Item(2, 6).Value = -0.05829021463284707

$ws.Cells.Item(3, 4).Value = -0.020236691044003
$ws.Cells.Item(3, 5).Value = -0.2736187788046003
$ws.Cells.Item(3, 6).Value = [double]"-3.243416007105181E-05"

$ws.Cells.Item(4, 1).Value = $null
$ws.Cells.Item(4, 4).Value = -0.04202283529381859
$ws.Cells.Item(4, 5).Value = -0.04602167871213612
$ws.Cells.Item(4, 6).Value = -0.05663840587470741

$ws.Cells.Item(5, 4).Value = -0.02075516704708091
$ws.Cells.Item(5, 5).Value = -0.2789229934737962
$ws.Cells.Item(5, 6).Value = [double]"-3.298510209777308E-05"

$ws.Cells.Item(6, 1).Value = $null
$ws.Cells.Item(6, 4).Value = -0.04309948520826185
$ws.Cells.Item(6, 5).Value = -0.04691382823634792
$ws.Cells.Item(6, 6).Value = -0.05760049270089667

$ws.Cells.Item(7, 4).Value = -0.025777836270233
$ws.Cells.Item(7, 5).Value = -0.4024328803204771
$ws.Cells.Item(7, 6).Value = [double]"-4.745086983336472E-05"

$ws.Cells.Item(8, 1).Value = $null
$ws.Cells.Item(8, 4).Value = -0.0535293920068044
$ws.Cells.Item(8, 5).Value = -0.06768773986282094
$ws.Cells.Item(8, 6).Value = -0.08286145282759177

$ws.Cells.Item(9, 4).Value = -0.02692763198875292
$ws.Cells.Item(9, 5).Value = -0.3758758581391763
$ws.Cells.Item(9, 6).Value = [double]"-4.438069385864546E-05"

$ws.Cells.Item(10, 4).Value = -0.05591701931187308
$ws.Cells.Item(10, 5).Value = -0.06322094577902851
$ws.Cells.Item(10, 6).Value = -0.07750013400256306

$ws.Cells.Item(11, 1).Value = "G_HepTh"
$ws.Cells.Item(11, 4).Value = -0.02792797929411883
$ws.Cells.Item(11, 5).Value = -0.3806207792499032
$ws.Cells.Item(11, 6).Value = [double]"-4.492923965179361E-05"

$ws.Cells.Item(12, 4).Value = -0.05799430704426965
$ws.Cells.Item(12, 5).Value = -0.06401902417052754
$ws.Cells.Item(12, 6).Value = -0.07845803638712082

$ws.Cells.Item(13, 1).Value = "G_HepTh"
$ws.Cells.Item(13, 4).Value = -0.02641307679014976
$ws.Cells.Item(13, 5).Value = -0.486126745848356
$ws.Cells.Item(13, 6).Value = [double]"-5.652978684547987E-05"

$ws.Cells.Item(14, 4).Value = -0.05484851120877158
$ws.Cells.Item(14, 5).Value = -0.08176474220282265
$ws.Cells.Item(14, 6).Value = -0.09871558271745171

$ws.Cells.Item(15, 1).Value = "G_HepTh"
$ws.Cells.Item(15, 4).Value = -0.02809375684972226
$ws.Cells.Item(15, 5).Value = -0.4651478862989276
$ws.Cells.Item(15, 6).Value = [double]"-5.434163712310543E-05"

$ws.Cells.Item(16, 4).Value = -0.0583385551676109
$ws.Cells.Item(16, 5).Value = -0.07823617468947844
$ws.Cells.Item(16, 6).Value = -0.09489450913888584

$ws.Cells.Item(17, 1).Value = "G_HepTh"
$ws.Cells.Item(17, 4).Value = -0.02878207076648515
$ws.Cells.Item(17, 5).Value = -0.4283807008613367
$ws.Cells.Item(17, 6).Value = [double]"-5.014859213172426E-05"

$ws.Cells.Item(18, 4).Value = -0.05976788480908614
$ws.Cells.Item(18, 5).Value = -0.07205206845689158
$ws.Cells.Item(18, 6).Value = -0.08757237150521491

